$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (Giorno_Cont, Giorno(date text), Infetti, Decessi, Ricoverati)
$data = @(
    @(118, "5/19/20", 232037, 27778, 150376),
    @(119, "5/20/20", 232555, 27888, 150376),
    @(120, "5/21/20", 233037, 27940, 150376),
    @(121, "5/22/20", 234824, 28628, 150376),
    @(122, "5/23/20", 235290, 28678, 150376),
    @(123, "5/24/20", 235772, 28752, 150376),
    @(124, "5/25/20", 235400, 26834, 150376),
    @(125, "5/26/20", 236259, 27117, 150376),
    @(126, "5/27/20", 236259, 27117, 150376),
    @(127, "5/28/20", 237906, 27119, 150376),
    @(128, "5/29/20", 238564, 27121, 150376),
    @(129, "5/30/20", 239228, 27125, 150376),
    @(130, "5/31/20", 239479, 27127, 150376),
    @(131, "6/1/20", 239638, 27127, 150376)
)

$startRow = 120
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rowValues[0]

    # Prefix the date-like string with an apostrophe so Excel stores it as
    # literal text (shared string) instead of auto-converting it to a date
    # serial number, then clear the cell formatting so no number-format
    # style is left attached to the cell itself.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.Value = "'" + $rowValues[1]
    $dateCell.ClearFormats()

    $ws.Cells.Item($row, 3).Value = $rowValues[2]
    $ws.Cells.Item($row, 4).Value = $rowValues[3]
    $ws.Cells.Item($row, 5).Value = $rowValues[4]
}
